$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Add new header cells (E1, F1) and new data cells (E2, F2)
# (order matters for shared-string table ordering)
$ws.Range("E1").Value = "Error message"
$ws.Range("E2").Value = "Please complete all required fields before sending."
$ws.Range("F1").Value = "Sign in txt"
$ws.Range("F2").Value = "Thank you"

# Adjust column E width (target stored width 13.6328125 characters)
$ws.Columns.Item(5).ColumnWidth = 12.8333333333333

# Configure page setup (paper size + orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the active selection to E2
$ws.Range("E2").Select()
